# DP: change output_inputs ad_hoc file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the "software location" description and path to point at the
# new basic-forecast tool output (was the ad_hoc tool folder).
$ws.Range("A2").Value = "מיקום תוכנת תחזית בסיס"
$ws.Range("B2").Value = "W:\Data\Forecast\Tools\forecast_git\create_forecast_basic\current"

# Row 6 was blank; fill it in with a new input/output pair describing where
# the base-forecast output (by version) lives.
$ws.Range("A6:B6").Style = "Normal"
$ws.Range("A6").Value = "מיקום פלט תחזית בסיס לפי גירסא"
$ws.Range("B6").Value = "W:\Data\Forecast\forecast_by_version\V4\BASE_YEAR"
$ws.Rows.Item(6).AutoFit()

# Match the saved selection state (A1:B6 highlighted).
$ws.Range("A1:B6").Select() | Out-Null
